# "Generate Report for Handback" -- reflect that the 17b0f83a... and
# a89c913d... files have now been handed back (in sync with en-US) for
# both the zh-cn and de-de locales: update each locale's Status column,
# fill in the "Latest Target File" / "Latest Handback File" columns (with
# hyperlinks back to the source .md on GitHub) and the "Latest Handback
# DateTime" column, then refresh the Overview roll-up and widen the
# columns that now hold the longer text.

$wb = $excel.ActiveWorkbook

$mdBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/de8c7e1762b83471abb4acb3db5101c7ec5b8ba5/e2e/"
$status = "Handed back: in sync with en-US"

$file1 = "17b0f83a-ffac-4a0d-a1c9-1fdcd4ba2337.md"
$file2 = "a89c913d-7275-4b74-9830-1cf7bfad5cf0.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $status
$wsZh.Range("C3").Value = $status

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($mdBase + $file1), $null, $null, $file1)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($mdBase + $file2), $null, $null, $file2)

$wsZh.Range("J2").Value = "17b0f83a-ffac-4a0d-a1c9-1fdcd4ba2337.8b99ba26d7abfa6a24388a680e1c3d5b27be8a26.zh-cn.xlf"
$wsZh.Range("J3").Value = "a89c913d-7275-4b74-9830-1cf7bfad5cf0.69bf9ec860a7de7189a950db8e93abd6940a61a5.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-13 16:33:29"
$wsZh.Range("K3").Value = "2016-08-13 16:33:29"

$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(9).ColumnWidth = 39.15
$wsZh.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $status
$wsDe.Range("C3").Value = $status

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($mdBase + $file1), $null, $null, $file1)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($mdBase + $file2), $null, $null, $file2)

$wsDe.Range("J2").Value = "17b0f83a-ffac-4a0d-a1c9-1fdcd4ba2337.8b99ba26d7abfa6a24388a680e1c3d5b27be8a26.de-de.xlf"
$wsDe.Range("J3").Value = "a89c913d-7275-4b74-9830-1cf7bfad5cf0.69bf9ec860a7de7189a950db8e93abd6940a61a5.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-13 16:33:39"
$wsDe.Range("K3").Value = "2016-08-13 16:33:39"

$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(9).ColumnWidth = 39.15
$wsDe.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# Overview roll-up sheet -- same status text shows in the per-locale
# summary columns, so widen them to match.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $status
$wsOverview.Range("F2").Value = $status
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status

$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1
